# Apply the edits described by the commit diff:
#  1. Shared string "Mean peak period" (A11 on Sheet1) -> "Mean wave period"
#  2. Sheet1's active selection moves from E10 to A13
#  3. Row heights for rows 9 and 10 shrink from 72 to 29.25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Update the label text in A11 (shared string "Mean peak period" -> "Mean wave period")
$ws.Range("A11").Value = "Mean wave period"

# 2. Adjust row heights for rows 9 and 10
$ws.Rows.Item(9).RowHeight = 29.25
$ws.Rows.Item(10).RowHeight = 29.25

# 3. Move the active selection to A13
$ws.Activate()
$ws.Range("A13").Select()
